$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header text for column B
$ws.Range("B1").Value = "TOOP Document Type Identifier"

# Remove the now-empty trailing rows (4-7) that only contained blank formatted cells
$ws.Rows("4:7").Delete()

# The wrapping text formatting is no longer needed on the remaining rows
$ws.Range("A1:D3").WrapText = $false

# Row 3 no longer needs the extra wrap height now that wrapping is gone;
# let Excel recompute the natural (non-custom) row height
$ws.Rows("1:3").AutoFit()

# Update the active selection to match the saved state
$ws.Range("D3").Select()
